# Card13: add a new "Serviced by " column (O) after the "Event" column (N),
# and tidy up the "Event " header (drop its trailing space).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")
$ws.Activate()

# --- N1: "Event " -> "Event" (trailing space removed) ---
$ws.Range("N1").Value = "Event"

# --- O1: new header "Serviced by ", styled like the rest of row 1 ---
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)   # xlPasteFormats: inherit the bold/centered header style
$ws.Range("O1").Value = "Serviced by "

# --- N2:N13 were blank placeholders; the edit fills them with "nan" ---
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}

# --- O2:O13: new column's data rows are blank (empty string) cells ---
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 15).Value = "'"       # write as an empty text entry
    $ws.Cells.Item($r, 15).Style = "Normal"  # drop the quote-prefix format it implies
}

Write-Output "Card13: added 'Serviced by ' column (O) and fixed 'Event' header"
